$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.148.36'
$ws.Range("E2").Value = '  +11.74%  '
$ws.Range("D3").Value = '3.278.28'
$ws.Range("E3").Value = '  +7.20%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '398.92'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.72'
$ws.Range("E6").Value = '  +10.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +5.21%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").Value = '  +7.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.71'
$ws.Range("E10").Value = '  +8.31%  '
$ws.Range("E11").Value = '  +12.34%  '
$ws.Range("E12").Value = '  +2.63%  '
$ws.Range("D13").Value = '3.780.10'
$ws.Range("E13").Value = '  +6.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.26'
$ws.Range("E14").Value = '  +5.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.13'
$ws.Range("E15").Value = '  +6.50%  '
$ws.Range("D16").Value = '3.267.33'
$ws.Range("E16").Value = '  +7.18%  '
$ws.Range("E17").Value = '  +4.38%  '
$ws.Range("E18").Value = '  +3.57%  '
$ws.Range("D19").Value = '56.955.13'
$ws.Range("E19").Value = '  +11.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.34'
$ws.Range("E20").Value = '  +6.18%  '
$ws.Range("E21").Value = '  +9.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.05'
$ws.Range("E22").Value = '  +6.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '303.45'
$ws.Range("E23").Value = '  +15.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.32'
$ws.Range("E24").Value = '  +8.25%  '
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.14'
$ws.Range("E26").Value = '  +3.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.42'
$ws.Range("E27").Value = '  +6.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.38'
$ws.Range("E28").Value = '  +5.41%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.171'
$ws.Range("E29").Value = '  +5.18%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.32'
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +6.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.05'
$ws.Range("E33").Value = '  +4.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '38.68'
$ws.Range("E34").Value = '  +8.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0488'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.15'
$ws.Range("E36").Value = '  +4.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.88'
$ws.Range("E37").Value = '  +3.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("E38").Value = '  +30.36%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.52'
$ws.Range("E39").Value = '  +5.81%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.77'
$ws.Range("E41").Value = '  +7.20%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.95'
$ws.Range("E42").Value = '  +6.71%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '133.91'
$ws.Range("E43").Value = '  +3.98%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.00'
$ws.Range("E44").Value = '  +6.84%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.120'
$ws.Range("E45").Value = '  +4.64%  '
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.28'
$ws.Range("E47").Value = '  +2.79%  '
$ws.Range("D48").Value = '2.158.70'
$ws.Range("E48").Value = '  +4.62%  '
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.02'
$ws.Range("E51").Value = '  +42.52%  '
